# Remove the "References" heading and its bulleted reference-list entries
# from the end of the document, leaving the preceding paragraph (ending in
# "...the increase of the AUC can be reported.") directly followed by the
# trailing empty paragraph.

$d = $word.ActiveDocument

# Locate the "References" heading paragraph and the last reference-list
# paragraph (the one ending in "...2750-2758") by scanning from the end of
# the document, so the script is resilient to any paragraph-count drift
# elsewhere in the document.
$paragraphs = $d.Paragraphs
$count = $paragraphs.Count

$startIndex = 0
$endIndex = 0

for ($i = $count; $i -ge 1; $i--) {
    $text = $paragraphs.Item($i).Range.Text
    if ($endIndex -eq 0 -and $text -like "*2750-2758*") {
        $endIndex = $i
    }
    if ($text.Trim() -eq "References") {
        $startIndex = $i
        break
    }
}

if ($startIndex -gt 0 -and $endIndex -ge $startIndex) {
    $start = $paragraphs.Item($startIndex).Range.Start
    $end = $paragraphs.Item($endIndex).Range.End
    $r = $d.Range($start, $end)
    $r.Delete()
    Write-Output "Removed paragraphs $startIndex..$endIndex"
} else {
    Write-Output "References section not found (start=$startIndex end=$endIndex)"
}
